$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) ---
$ws.Range("A1").Value = "Include"
$ws.Range("B1").Value = "Criteria"
$ws.Range("C1").Value = "Exceptional"
$ws.Range("D1").Value = "Proficient"
$ws.Range("E1").Value = "Basic"
$ws.Range("F1").Value = "Limited"
$ws.Range("G1").Value = "Missing"

# --- Row 2: Summary of Provided Service ---
$ws.Range("A2").Value = "y"
$ws.Range("B2").Value = "Summary of Provided Service"
$ws.Range("C2").Value = "It is obvious from the presentation that your team is the best group to provide exactly the service you are proposing to most effectively solve your client's business problem. "
$ws.Range("D2").Value = "The product/service offering is condensed to one clear and crisp statement that does not leave any remaining questions on what will be offered to the stakeholders "
$ws.Range("E2").Value = "The services/products to be provided to the stakeholder(s) are outlined. But the description could have been made more compact and focused and more clearly formulated. After the description there are still some questions on what your company is exactly going to offer as primary service/product "
$ws.Range("F2").Value = "From the top statement on the slide it is not clear what the excact service / product is that is being provided to the stakeholder."
$ws.Range("G2").Value = "This item is missing."

# --- Row 3: Description of the Problem/Opportunity ---
$ws.Range("A3").Value = "y"
$ws.Range("B3").Value = "Description of the Problem/Opportunity"
$ws.Range("C3").Value = "A convincing story is created for the audience around the actual business problem(s) at hand and where the opportunities lie. This is communicated clearly and logically in a few concise statements."
$ws.Range("D3").Value = "The customer and their pain points or the missing opportunities they currently are not taking advantage of are described in a few sentences and a clear picture emerges to the audience of what the business problem/opportunity is."
$ws.Range("E3").Value = "The customer's main pain points are described but there are still questions remaining or too much time and detail is spent on describing instead of connecting them with the suggested solution."
$ws.Range("F3").Value = "The customer's pain points are not or not clearly addressed."
$ws.Range("G3").Value = "This item is missing."

# --- Row 4: Explanation of the Approach/Solution ---
$ws.Range("A4").Value = "y"
$ws.Range("B4").Value = "Explanation of the Approach/Solution"
$ws.Range("C4").Value = "A well-crafted story emerges from your presentation of how your solution will address the customer's business problem in a language that a non-technical stakeholder will understand. Tone, pace and delivery are well adjusted to communicate the pitch effectively in the 1 min time limit allotted. "
$ws.Range("D4").Value = "An explanation is provided of how the service/product/solution you are proposing will solve for the problem. Points flow logically. Usage scenarios are described and how will you address them in simple, customer-understood language. "
$ws.Range("E4").Value = "The description of how your solution fits the businesx problem is logically sound but misses emphasis of key aspects and is not able to draw in the audience to be excited or encouraged about supporting your proposal and following through on the action(s) requested of them."
$ws.Range("F4").Value = "There are logical gaps in how your solution connects with the business problem of your customer and the description in the presentation is not easy to follow."
$ws.Range("G4").Value = "This item is missing."

# --- Row 5: Description of the Customer Experience and Capturing the Business Value Provided ---
$ws.Range("A5").Value = "y"
$ws.Range("B5").Value = "Description of the Customer Experience and Capturing the Business Value Provided"
$ws.Range("C5").Value = "In the presentation it can be clearly seen that the pitch is conducted with customer's perspective in mind. The solution you are presenting is clearly solving their problem, or an aspect of it, and the audience can convincingly see the business value in your solution."
$ws.Range("D5").Value = "A customer journey is painted from the perspective of the customer where it is clear from the description how will discover and use your solution. You can outline your value propsition, i.e. how the customer will react / feel about your solution."
$ws.Range("E5").Value = "The business problem and your team's solution are connected and this can be seen from your presentation. But the asrgumentation does not flow easily, key questions remain whether a better aproach exists or whether your solution would truly create business value for the stakeholder."
$ws.Range("F5").Value = "There is a disconnect between the customer's problem(s) and the solution your team is proposing. There is the impression that the presented problem serves to justify your team's solution and not the opposite. The business value of your solution is not convincingly presented to the audience."
$ws.Range("G5").Value = "This item is missing."

# --- Row 6: Call to Action ---
$ws.Range("A6").Value = "y"
$ws.Range("B6").Value = "Call to Action"
$ws.Range("C6").Value = "A call to action very nicrly flows from the flow of the presentation based on the action items listed above."
$ws.Range("D6").Value = "A call to action is clearly formulated and to the point and can be directly derived from the explanations on the above items. There are no questions remaining of what is requested from the stakeholders."
$ws.Range("E6").Value = "A call to action is formulated but could be made more precise and some questions remain."
$ws.Range("F6").Value = "A call to action for the stakeholders is not clearly formulated and confusing."
$ws.Range("G6").Value = "This item is missing."

# --- Row 7: Powerpoint Slide Layout ---
$ws.Range("A7").Value = "y"
$ws.Range("B7").Value = "Powerpoint Slide Layout"
$ws.Range("C7").Value = "The single slide for this solution pitch is exceptionally well designed."
$ws.Range("D7").Value = "The slide layout shows a professional layout and presents all the information concisely with a straightforward logical flow."
$ws.Range("E7").Value = "The slide layout is functional in structure and relating the different apoint above but could be made more effective and further improved to communicate the key messages."
$ws.Range("F7").Value = "The slide layout appears disorganized and hard to follow. A logical flow of information from top to bottom is missing."
$ws.Range("G7").Value = "This item is missing."

# --- Row heights (auto-fit heights that Excel computed for the wrapped text) ---
$ws.Rows.Item(2).RowHeight = 119
$ws.Rows.Item(3).RowHeight = 85
$ws.Rows.Item(4).RowHeight = 119
$ws.Rows.Item(5).RowHeight = 119
$ws.Rows.Item(6).RowHeight = 85
$ws.Rows.Item(7).RowHeight = 68

# --- Selection ---
[void]$ws.Range("E4").Select()
